# "Got 2nd Alt working"
#
# 1. Rename the "Alt1 (2)" sheet to "Alt2"
# 2. Replace its B11:B15 array formula with a BYROW/LAMBDA wrapper around
#    REGEXEXTRACT (so it extracts per-row instead of spilling the whole
#    range at once), using a non-greedy pattern with both lookbehind and
#    lookahead.
# 3. Move the sheet's active selection from B16 to B14.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Alt1 (2)")
$ws.Name = "Alt2"

$rng = $ws.Range("B11:B15")
$rng.FormulaArray = '=_xlfn.BYROW(B3:B7,_xlfn.LAMBDA(x,_xlfn.REGEXEXTRACT(x,"(?<=\().*?(?=\))|(?<=\[).*?(?=\])|(?<=\{).*?(?=\})")))'

$ws.Activate() | Out-Null
$ws.Range("B14").Select() | Out-Null
